$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.292064567892659 / 100000
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("G2").Value = 2196001.309727008

# Row 3
$ws.Range("B3").Value = 0.2917716402565462
$ws.Range("C3").Value = 2.903117352495289 * 100000000000000000
$ws.Range("D3").Value = 19477208507.93593
$ws.Range("E3").Value = 91228006295.30009
$ws.Range("G3").Value = 2.903118459547437 * 100000000000000000
